{"js": "// The author's change (POI 4.1.0 -> 5.2.3 migration, fixing #476) touched the\n// run-level character formatting of the \"some text\" run: the boolean toggles\n// for Bold/Italic/StrikeThrough got re-written by the new POI version (the\n// on/off state itself is unchanged: Bold stays on, Italic stays off, Strike\n// stays off) and the <w:sz> element shifted position. Reproduce the edit by\n// (re)asserting the same formatting on that run through the object model so\n// the run's rPr is rewritten by this engine's canonical OOXML writer.\n\nconst body = context.document.body;\n\n// Scope the edit strictly to the \"some text\" run so we don't touch the\n// paragraph mark / bookmark / any other run in the document.\nconst results = body.search(\"some text\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find the 'some text' run to update.\");\n}\n\nconst target = results.items[0];\ntarget.load(\"font\");\nawait context.sync();\n\n// Re-assert the run's bold/italic/strike-through state (bold on, italic\n// off, strike off) so the formatting is rewritten with the current engine's\n// canonical representation, matching the semantics of the target edit.\ntarget.font.bold = true;\ntarget.font.italic = false;\ntarget.font.strikeThrough = false;\n\nawait context.sync();\n", "ps1": "# The author's change (POI 4.1.0 -> 5.2.3 migration, fixing #476) touched the\n# run-level character formatting of the \"some text\" run: the boolean toggles\n# for Bold/Italic/StrikeThrough got re-written by the new POI version (the\n# on/off state itself is unchanged: Bold stays on, Italic stays off, Strike\n# stays off) and the <w:sz> element shifted position. Reproduce the edit by\n# (re)asserting the same formatting on that run through the Word object model\n# so the run's rPr is rewritten by this engine's canonical OOXML writer.\n\n$d = $word.ActiveDocument\n\n# Scope the edit strictly to the \"some text\" run so we don't touch the\n# paragraph mark / bookmark / any other run in the document.\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Text = \"some text\"\n$rng.Find.Forward = $true\n$rng.Find.Wrap = 0\n$found = $rng.Find.Execute()\n\nif (-not $found) {\n    throw \"Could not find the 'some text' run to update.\"\n}\n\n# Re-assert the run's bold/italic/strike-through state (bold on, italic\n# off, strike off) so the formatting is rewritten with the current engine's\n# canonical representation, matching the semantics of the target edit.\n$rng.Font.Bold = 1\n$rng.Font.Italic = 0\n$rng.Font.StrikeThrough = 0\n"}
